$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 316-317; existing rows 316-403 shift down to 318-405,
# matching the new dimension A1:R405.
$ws.Rows("316:317").Insert()

# Populate the two newly-inserted rows with this week's price data
# (row 316 = "Primera" quality, row 317 = "Segunda" quality), keeping the
# other descriptive columns identical to the block they precede.

# Row 316
$ws.Cells.Item(316, 1).Value = 8
$ws.Cells.Item(316, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(316, 3).Value = "Coquimbo"
$ws.Cells.Item(316, 4).Value = 44964
$ws.Cells.Item(316, 5).Value = 4
$ws.Cells.Item(316, 6).Value = 100114014
$ws.Cells.Item(316, 7).Value = "Betarraga"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 1780
$ws.Cells.Item(316, 11).Value = 500
$ws.Cells.Item(316, 12).Value = 600
$ws.Cells.Item(316, 13).Value = 550
$ws.Cells.Item(316, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(316, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(316, 16).Value = 183
$ws.Cells.Item(316, 17).Value = 3
$ws.Cells.Item(316, 18).Value = "Hortaliza"

# Row 317
$ws.Cells.Item(317, 1).Value = 8
$ws.Cells.Item(317, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(317, 3).Value = "Coquimbo"
$ws.Cells.Item(317, 4).Value = 44964
$ws.Cells.Item(317, 5).Value = 4
$ws.Cells.Item(317, 6).Value = 100114014
$ws.Cells.Item(317, 7).Value = "Betarraga"
$ws.Cells.Item(317, 8).Value = "Sin especificar"
$ws.Cells.Item(317, 9).Value = "Segunda"
$ws.Cells.Item(317, 10).Value = 1360
$ws.Cells.Item(317, 11).Value = 400
$ws.Cells.Item(317, 12).Value = 450
$ws.Cells.Item(317, 13).Value = 425
$ws.Cells.Item(317, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(317, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(317, 16).Value = 142
$ws.Cells.Item(317, 17).Value = 3
$ws.Cells.Item(317, 18).Value = "Hortaliza"
